# Complete rebuild of database (db name: cmms2)
# Rename database references from "newcmms" to "cmms2" throughout the sheet,
# and move the active selection to A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 holds the database name, A4 holds the \c <dbname> command.
$ws.Range("A2").Value = "cmms2"
$ws.Range("A4").Value = "\c cmms2"

# Update the selected cell to A7 to match the saved view state.
$ws.Range("A7").Select()
